$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "44+38=82"
$t.Cell(1,2).Range.Text = "84-36=48"
$t.Cell(1,3).Range.Text = "70-13=57"
$t.Cell(1,4).Range.Text = "70-22=48"
$t.Cell(1,5).Range.Text = "75-16=59"
$t.Cell(2,1).Range.Text = "60-3=57"
$t.Cell(2,2).Range.Text = "36+59=95"
$t.Cell(2,3).Range.Text = "25+7=32"
$t.Cell(2,4).Range.Text = "3+49=52"
$t.Cell(2,5).Range.Text = "63-44=19"
$t.Cell(3,1).Range.Text = "17+55=72"
$t.Cell(3,2).Range.Text = "41-19=22"
$t.Cell(3,3).Range.Text = "25+49=74"
$t.Cell(3,4).Range.Text = "71-49=22"
$t.Cell(3,5).Range.Text = "9+54=63"
$t.Cell(4,1).Range.Text = "24+17=41"
$t.Cell(4,2).Range.Text = "92-24=68"
$t.Cell(4,3).Range.Text = "52-9=43"
$t.Cell(4,4).Range.Text = "49+38=87"
$t.Cell(4,5).Range.Text = "18+37=55"
$t.Cell(5,1).Range.Text = "71-64=7"
$t.Cell(5,2).Range.Text = "91-87=4"
$t.Cell(5,3).Range.Text = "5+66=71"
$t.Cell(5,4).Range.Text = "68+9=77"
$t.Cell(5,5).Range.Text = "53-28=25"
$t.Cell(6,1).Range.Text = "62-56=6"
$t.Cell(6,2).Range.Text = "41-38=3"
$t.Cell(6,3).Range.Text = "34+57=91"
$t.Cell(6,4).Range.Text = "84-29=55"
$t.Cell(6,5).Range.Text = "58+3=61"
$t.Cell(7,1).Range.Text = "50-14=36"
$t.Cell(7,2).Range.Text = "46+6=52"
$t.Cell(7,3).Range.Text = "65+29=94"
$t.Cell(7,4).Range.Text = "58+27=85"
$t.Cell(7,5).Range.Text = "34+9=43"
$t.Cell(8,1).Range.Text = "10-6=4"
$t.Cell(8,2).Range.Text = "16+58=74"
$t.Cell(8,3).Range.Text = "28+57=85"
$t.Cell(8,4).Range.Text = "46+39=85"
$t.Cell(8,5).Range.Text = "63-6=57"
$t.Cell(9,1).Range.Text = "72-28=44"
$t.Cell(9,2).Range.Text = "28+33=61"
$t.Cell(9,3).Range.Text = "35+8=43"
$t.Cell(9,4).Range.Text = "7+48=55"
$t.Cell(9,5).Range.Text = "49+27=76"
$t.Cell(10,1).Range.Text = "63-37=26"
$t.Cell(10,2).Range.Text = "18+78=96"
$t.Cell(10,3).Range.Text = "55+6=61"
$t.Cell(10,4).Range.Text = "54-5=49"
$t.Cell(10,5).Range.Text = "80-6=74"
$t.Cell(11,1).Range.Text = "48+34=82"
$t.Cell(11,2).Range.Text = "84-39=45"
$t.Cell(11,3).Range.Text = "67+7=74"
$t.Cell(11,4).Range.Text = "94-68=26"
$t.Cell(11,5).Range.Text = "60-2=58"
$t.Cell(12,1).Range.Text = "35+56=91"
$t.Cell(12,2).Range.Text = "74-37=37"
$t.Cell(12,3).Range.Text = "8+75=83"
$t.Cell(12,4).Range.Text = "66-17=49"
$t.Cell(12,5).Range.Text = "29+15=44"
$t.Cell(13,1).Range.Text = "16+46=62"
$t.Cell(13,2).Range.Text = "21-19=2"
$t.Cell(13,3).Range.Text = "28-19=9"
$t.Cell(13,4).Range.Text = "55+8=63"
$t.Cell(13,5).Range.Text = "39+4=43"
$t.Cell(14,1).Range.Text = "63-19=44"
$t.Cell(14,2).Range.Text = "84-66=18"
$t.Cell(14,3).Range.Text = "34+19=53"
$t.Cell(14,4).Range.Text = "85-47=38"
$t.Cell(14,5).Range.Text = "61-49=12"
$t.Cell(15,1).Range.Text = "92-76=16"
$t.Cell(15,2).Range.Text = "92-55=37"
$t.Cell(15,3).Range.Text = "33+59=92"
$t.Cell(15,4).Range.Text = "36+49=85"
$t.Cell(15,5).Range.Text = "82-38=44"
$t.Cell(16,1).Range.Text = "90-38=52"
$t.Cell(16,2).Range.Text = "9+19=28"
$t.Cell(16,3).Range.Text = "75-46=29"
$t.Cell(16,4).Range.Text = "55+36=91"
$t.Cell(16,5).Range.Text = "60-26=34"
$t.Cell(17,1).Range.Text = "15+6=21"
$t.Cell(17,2).Range.Text = "40-7=33"
$t.Cell(17,3).Range.Text = "75-9=66"
$t.Cell(17,4).Range.Text = "27+8=35"
$t.Cell(17,5).Range.Text = "18+9=27"
$t.Cell(18,1).Range.Text = "57+16=73"
$t.Cell(18,2).Range.Text = "45+16=61"
$t.Cell(18,3).Range.Text = "80-11=69"
$t.Cell(18,4).Range.Text = "38-29=9"
$t.Cell(18,5).Range.Text = "73-64=9"
$t.Cell(19,1).Range.Text = "49+18=67"
$t.Cell(19,2).Range.Text = "44-6=38"
$t.Cell(19,3).Range.Text = "66+8=74"
$t.Cell(19,4).Range.Text = "90-4=86"
$t.Cell(19,5).Range.Text = "28+15=43"
$t.Cell(20,1).Range.Text = "10-9=1"
$t.Cell(20,2).Range.Text = "55+16=71"
$t.Cell(20,3).Range.Text = "90-74=16"
$t.Cell(20,4).Range.Text = "76+19=95"
$t.Cell(20,5).Range.Text = "17+57=74"
